# Rework the scraper output sheet: the price and item-name/description
# columns are swapped (what used to be column B -- text -- becomes column
# C, and what used to be column C -- the price number -- becomes column
# B), headers are relabelled, and the sheet view / column widths are
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Capture the current (pre-edit) values for columns B and C, rows
#    2-7, before anything gets overwritten.
# ---------------------------------------------------------------------
$oldB = @{}
$oldC = @{}
for ($r = 2; $r -le 7; $r++) {
    $oldB[$r] = $ws.Cells.Item($r, 2).Value2
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
}

# ---------------------------------------------------------------------
# 2. Swap the cell formatting between columns B and C (rows 2-7) so the
#    numeric-price look (currency number format) ends up on column B and
#    the text look ends up on column C. Stash column C's original format
#    in an unused helper column first since PasteSpecial overwrites the
#    destination in place.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Copy() | Out-Null
    $ws.Cells.Item($r, 26).PasteSpecial($xlPasteFormats) | Out-Null
}
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 2).Copy() | Out-Null
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats) | Out-Null
}
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 26).Copy() | Out-Null
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats) | Out-Null
}
$ws.Range("Z2:Z7").Clear() | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Write the swapped values: B gets the old price number (was in C),
#    C gets the old text (was in B).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 2).Value = [double]$oldC[$r]
    $ws.Cells.Item($r, 3).Value = $oldB[$r]
}

# ---------------------------------------------------------------------
# 4. Header row: A1 stays "sku", B1 becomes "our_price", C1 becomes
#    "Item Name" and picks up column A's header-ish look.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 3).PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "our_price"
$ws.Range("C1").Value = "Item Name"

# ---------------------------------------------------------------------
# 5. Column widths: B narrows to fit the numeric price (matches the old
#    column C width, ~11.14 chars), C drops its bestFit override and
#    falls back to the sheet default width (~9.14 chars). The inputs
#    below are tuned to the engine's column-width rounding so the saved
#    width lands as close as possible to those targets.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.33
$ws.Columns.Item(3).ColumnWidth = 8.33

# ---------------------------------------------------------------------
# 6. Selection moves to E8.
# ---------------------------------------------------------------------
$ws.Range("E8").Select() | Out-Null
